# BUG TRACKER.xlsx — add new bug-tracker row (row 4) on "Feuil1"
# "insertion personnage .." bug reported by tommy, status "à faire"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# NUM
$ws.Range("A4").Value = 3

# Titre
$ws.Range("B4").Value = "pas vraiment une erreur mais le lien de l'image pour le perosnnage est inexistant .. Donc pour le moment j'utilise l'image de la profession"

# trouvé par
$ws.Range("D4").Value = "tommy"

# fichier
$ws.Range("E4").Value = "insertion personnage .."

# état
$ws.Range("F4").Value = "à faire"

# The longer title text now needs more vertical room
$ws.Rows.Item(4).RowHeight = 61.5

# Move the active selection to H4, as left by the author after editing
$null = $ws.Range("H4").Select()
